$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("seats")

$ws.Range("A1").Value = 2
$ws.Range("B1").Value = "klf"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "klf"

$ws.Range("A3").Value = 25
$ws.Range("B3").Value = "admin"

$ws.Range("A4").Value = 23
$ws.Range("B4").Value = "admin"
